# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns populated, with a new hyperlink on the
#    "Latest Target File" cell
#  - Overview + per-locale sheets widen a few columns to fit the new content

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across Overview!E2/F2 and the locale sheets' Status column C2)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet row 2: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$zhcn.Range("I2").Value = "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md", "", "", "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md") | Out-Null
$zhcn.Range("J2").Value = "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.f4f378ba79d76271ad4629a62a122341b2ecc7a1.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-17 10:55:40"

# ---------------------------------------------------------------------------
# 3. de-de sheet row 2: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$dede.Range("I2").Value = "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md", "", "", "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.md") | Out-Null
$dede.Range("J2").Value = "81ff3386-6d10-4c04-ad6f-359ca4ef08bd.f4f378ba79d76271ad4629a62a122341b2ecc7a1.de-de.xlf"
$dede.Range("K2").Value = "2016-08-17 10:55:48"

# ---------------------------------------------------------------------------
# 4. Column widths - widen to fit the newly populated long filenames
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
